{"js": "// Remove the first paragraph of the document body\n// (\"Last call for the course on Text Mining with R, held next week in\n// Leuven, Belgium on April 1-2. ...\") leaving the rest of the body intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.delete();\nawait context.sync();\n", "ps1": "# Remove the first paragraph of the document body\n# (\"Last call for the course on Text Mining with R, held next week in\n# Leuven, Belgium on April 1-2. ...\") leaving the rest of the body intact.\n$d = $word.ActiveDocument\n$d.Paragraphs(1).Range.Delete()\n"}
